$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '29.557.08'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value2 = '  +0.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '1.923.22'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value2 = '  +0.53%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.011'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value2 = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '326.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = '  +0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '1.010'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = '  +0.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.4820'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value2 = '  -0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.4071'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = '  +0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.08232'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = '  +0.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '1.011'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = '  -0.15%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '23.65'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = '  +1.14%  '

$ws.Range("B12").Value2 = 'Polkadot'
$ws.Range("C12").Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '6.082'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = '  +1.34%  '

$ws.Range("B13").Value2 = 'WrappedEther'
$ws.Range("C13").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '1.894.64'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = '  -1.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '7.275'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = '  +2.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '91.74'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = '  +1.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '0.06878'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = '  +1.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '0.00001038'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = '  -0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '17.63'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = '  -0.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '1.009'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = '  +0.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '29.569.99'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = '  +0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '5.689'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = '  +1.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '11.92'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = '  +0.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.184'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = '  +0.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '2.155.40'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = '  +0.43%  '

$ws.Range("E26").Value2 = '  +0.23%  '

$ws.Range("E27").Value2 = '  +1.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '20.03'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = '  -0.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '2.100'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = '  +0.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '120.63'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = '  +0.82%  '

$ws.Range("E31").Value2 = '  -1.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '0.09640'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = '  +0.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '5.634'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = '  +1.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '3.550'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = '  -0.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '1.377'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value2 = '  -1.28%  '

$ws.Range("E36").Value2 = '  +4.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.02298'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = '  +1.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '1.190'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value2 = '  +0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '0.5957'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = '  +0.23%  '

$ws.Range("E40").Value2 = '  -0.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '7.897'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value2 = '  -0.59%  '

$ws.Range("E42").Value2 = '  -0.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '2.475'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = '  +0.59%  '

$ws.Range("B44").Value2 = 'WEMIXToken'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '1.279'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = '  -0.52%  '

$ws.Range("B45").Value2 = 'EnergySwap'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '12.45'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = '  +0.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '0.07497'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = '  -3.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '0.5573'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = '  +0.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '1.953'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = '  +0.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '119.21'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = '  +2.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '2.432'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = '  +3.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '72.22'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value2 = '  -0.69%  '
